$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Yadgir"
$ws.Range("G4").Value = "Yadgir"
$ws.Range("G5").Value = "Yadgir"
$ws.Range("G7").Value = "Yadgir"
$ws.Range("G8").Value = "Yadgir"

$ws.Range("G9").Value = "Ballari (Bellary)"
$ws.Range("G10").Value = "Ballari (Bellary)"
$ws.Range("G11").Value = "Ballari (Bellary)"
$ws.Range("G12").Value = "Ballari (Bellary)"
$ws.Range("G13").Value = "Ballari (Bellary)"
$ws.Range("G14").Value = "Ballari (Bellary)"
$ws.Range("G15").Value = "Ballari (Bellary)"
$ws.Range("G16").Value = "Ballari (Bellary)"
$ws.Range("G18").Value = "Ballari (Bellary)"
$ws.Range("G20").Value = "Ballari (Bellary)"
$ws.Range("G22").Value = "Ballari (Bellary)"
$ws.Range("G23").Value = "Ballari (Bellary)"
$ws.Range("G24").Value = "Ballari (Bellary)"
$ws.Range("G25").Value = "Ballari (Bellary)"
$ws.Range("G26").Value = "Ballari (Bellary)"
